$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 4.036825632626817
$ws.Range("C3").Value = -3.591129714716879
$ws.Range("E3").Value = -1.300150869559236
$ws.Range("C4").Value = 3.52827217675542
$ws.Range("E4").Value = 1.194245528719495
$ws.Range("C5").Value = 7.607887362976751
$ws.Range("E5").Value = 3.063639588842682
$ws.Range("C6").Value = 3.881953143326466
$ws.Range("E6").Value = 7.155859501153827
$ws.Range("C7").Value = 0.354775830825127
$ws.Range("E7").Value = 2.700497159199755
$ws.Range("C8").Value = 5.106323395421475
$ws.Range("E8").Value = 2.62739064366051
$ws.Range("C9").Value = 3.830515520137801
$ws.Range("E9").Value = 3.858663141671226
$ws.Range("C10").Value = 4.328608026086478
$ws.Range("E10").Value = 4.721407739775696
$ws.Range("C11").Value = 4.429102498614346
$ws.Range("E11").Value = 4.21218881008929
$ws.Range("C12").Value = 5.850954342715009
$ws.Range("E12").Value = 4.990046926794744
$ws.Range("C13").Value = 4.557673974453769
$ws.Range("E13").Value = 5.208344373007368
$ws.Range("C14").Value = 1.13158575217045
$ws.Range("E14").Value = 2.143123507515932
$ws.Range("C15").Value = -1.77012120409461
$ws.Range("E15").Value = -0.6252235182164778
$ws.Range("C16").Value = 1.286283684448075
$ws.Range("E16").Value = -2.183761975384579
$ws.Range("C17").Value = -0.7189954590872905
$ws.Range("E17").Value = 0.7347074005453758
$ws.Range("C18").Value = -0.382605475081077
$ws.Range("E18").Value = 0.3097078768351302
$ws.Range("C19").Value = 0.2094327661663842
$ws.Range("E19").Value = -0.6242159253788016
